# Actualización automática 2025-06-16 17:20:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D21").Value = 274.75
$ws1.Range("D22").Value = "1 de 20"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F21").Value = 274.75
$ws2.Range("F22").Value = -54.75999999999999

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 274.75
$ws3.Range("E3").Value = 3893.32156573679
$ws3.Range("F3").Value = 0.0659177741233031

$ws3.Range("D19").Value = -54.75999999999999
$ws3.Range("E19").Value = 50441.95762291769
$ws3.Range("F19").Value = -0.001086783996399383
